# Daily attendance processing - 2025-10-05 13:38:42
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Class Statistics block (K/L columns) ---
# L9/L10 are text percentages (not numeric) in the source sheet; a leading
# apostrophe forces Excel to store the literal text instead of auto-converting
# the percent-looking string into a numeric percentage.
$ws.Range("L6").Value = 49
$ws.Range("L7").Value = 7
$ws.Range("L9").Value = "'32.0%"
$ws.Range("L10").Value = "'48.6%"

# --- Recorded-by order shuffles (emails reordered, same sets) ---
$ws.Range("G8").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G14").Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G17").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"

# --- Group Statistics block (A4 row, row 18) ---
$ws.Range("O18").Value = 8
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = "'47.1%"
$ws.Range("S18").Value = "'51.8%"

$ws.Range("G25").Value = "Rania.a.youssef@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("G31").Value = "marian.samir@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G34").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
$ws.Range("G35").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

$ws.Range("H37").Value = "109/220"

$ws.Range("G45").Value = "System, backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G48").Value = "marian.samir@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G51").Value = "yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G52").Value = "Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"

# --- Row 54 (A4 / BIOCHEMISTRY LAB/CBL #1) went from Not Recorded -> Recorded ---
# Copy the formatting (fill/style) from an already-"Recorded" row (row 7) so the
# row picks up the same "Recorded" style index, then overwrite the actual values.
$ws.Range("A7:I7").Copy()
$ws.Range("A54:I54").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("G54").Value = "Kerelos.zareef@med.asu.edu.eg"
$ws.Range("H54").Value = "121/225"
$ws.Range("I54").Value = "Recorded"

$ws.Range("G62").Value = "System, backup@backdoor.com, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G65").Value = "marian.samir@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg"
$ws.Range("G68").Value = "yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G69").Value = "Salma.hassan@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G72").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G76").Value = "mariam.youssif.std@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G83").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"
$ws.Range("G85").Value = "yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

$ws.Range("G98").Value = "afaf.abdallah@med.asu.edu.eg, user@user.com, Walaa.h.ghanima@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
$ws.Range("G102").Value = "yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

$ws.Range("G116").Value = "enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G119").Value = "marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

$ws.Range("G133").Value = "enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G136").Value = "marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

$ws.Range("G150").Value = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"
$ws.Range("G151").Value = "marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg"
